$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

$ws.Range("D14").Value = "TIMESTAMP"
$ws.Range("D16").Value = "TIMESTAMP"

$ws.Range("D16").Select() | Out-Null
